$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc" (sheet1): update the day/night timeslice group strings ---
$wsUc = $wb.Worksheets.Item("ev_charging_uc")
$wsUc.Range("D11").Value = "S1aH3,S2aH2,S1aH2,S2aH3"
$wsUc.Range("D12").Value = "S2aH1,S1aH4,S1aH1,S2aH4"

# --- Sheet "ts12_clu" (sheet2): TFM_INS-AT hydro block updates ---
$wsTs = $wb.Worksheets.Item("ts12_clu")

# Commodity column AG (rows 11-18): "Elec" -> "ELC"
$wsTs.Range("AG11").Value = "ELC"
$wsTs.Range("AG12").Value = "ELC"
$wsTs.Range("AG13").Value = "ELC"
$wsTs.Range("AG14").Value = "ELC"
$wsTs.Range("AG15").Value = "ELC"
$wsTs.Range("AG16").Value = "ELC"
$wsTs.Range("AG17").Value = "ELC"
$wsTs.Range("AG18").Value = "ELC"

# AK/AL rows 11 and 12 swap (timeslice S1/S2 and matching ncap_afs values)
$wsTs.Range("AK11").Value = "S1"
$wsTs.Range("AL11").Value = 1.0373322535863025
$wsTs.Range("AK12").Value = "S2"
$wsTs.Range("AL12").Value = 0.16266774641369736

# Recalculate so dependent formula cells (e.g. H23/H24 HLOOKUP on sheet1) refresh
$excel.Calculate()
